$d = $word.ActiveDocument

$pairs = @(
    @("15×30=450", "19×65=1235"),
    @("42×35=1470", "12×15=180"),
    @("39×23=897", "96×99=9504"),
    @("48×20=960", "85×97=8245"),
    @("41×81=3321", "49×94=4606"),
    @("52×37=1924", "84×86=7224"),
    @("39×21=819", "88×91=8008"),
    @("92×16=1472", "58×71=4118"),
    @("57×85=4845", "19×87=1653"),
    @("91×69=6279", "83×88=7304"),
    @("38×24=912", "93×49=4557"),
    @("21×84=1764", "77×36=2772"),
    @("20×88=1760", "34×40=1360"),
    @("28×63=1764", "93×15=1395"),
    @("43×98=4214", "13×71=923"),
    @("79×25=1975", "76×53=4028"),
    @("45×60=2700", "54×87=4698"),
    @("42×27=1134", "16×46=736"),
    @("63×92=5796", "91×50=4550"),
    @("91×63=5733", "90×56=5040"),
    @("33×23=759", "76×63=4788"),
    @("15×24=360", "53×85=4505"),
    @("76×48=3648", "30×98=2940"),
    @("18×37=666", "11×56=616"),
    @("52×79=4108", "95×66=6270")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
